# Fixed SRS and product backlog issues
#
# The "Sprint 5" backlog had a numbering gap in column D (priority + story
# number): it jumped from H,6 straight to M,7 and from M,11 straight to
# L,12. This renumbers the whole sequence contiguously (H,1-4 / M,5-9 /
# L,10-22), which introduces two new "M" labels and two new "L" labels
# that previously did not exist anywhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 5")

$ws.Range("D8").Value  = "H, 1"
$ws.Range("D9").Value  = "H, 2"
$ws.Range("D10").Value = "H, 3"
$ws.Range("D11").Value = "H, 4"
$ws.Range("D12").Value = "M, 5"
$ws.Range("D13").Value = "M, 6"
$ws.Range("D14").Value = "M, 7"
$ws.Range("D15").Value = "M, 8"
$ws.Range("D16").Value = "M, 9"
$ws.Range("D17").Value = "L, 10"
$ws.Range("D18").Value = "L, 11"
$ws.Range("D19").Value = "L, 12"
$ws.Range("D20").Value = "L, 13"
$ws.Range("D21").Value = "L, 14"
$ws.Range("D22").Value = "L, 15"
$ws.Range("D23").Value = "L, 16"
$ws.Range("D24").Value = "L, 17"
$ws.Range("D25").Value = "L, 18"
$ws.Range("D26").Value = "L, 19"
$ws.Range("D27").Value = "L, 20"
$ws.Range("D28").Value = "L, 21"
$ws.Range("D29").Value = "L, 22"

# Restore/update the view state of the active sheet: scroll up a bit and
# move the selection to G27 (matches the saved workbookView after editing).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G27").Select()
